# Update "想去人数" (want-to-go count) figures in column F across the
# three affected worksheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1126
$ws1.Range("F4").Value  = 253
$ws1.Range("F5").Value  = 141
$ws1.Range("F6").Value  = 12100
$ws1.Range("F7").Value  = 48
$ws1.Range("F8").Value  = 89
$ws1.Range("F9").Value  = 11870
$ws1.Range("F10").Value = 4765
$ws1.Range("F11").Value = 580
$ws1.Range("F12").Value = 75
$ws1.Range("F13").Value = 26
$ws1.Range("F15").Value = 87
$ws1.Range("F16").Value = 934
$ws1.Range("F19").Value = 59
$ws1.Range("F20").Value = 5216

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 5

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1126
$ws4.Range("F4").Value  = 253
$ws4.Range("F5").Value  = 141
$ws4.Range("F6").Value  = 5
$ws4.Range("F8").Value  = 12100
$ws4.Range("F9").Value  = 48
$ws4.Range("F10").Value = 89
$ws4.Range("F11").Value = 11870
$ws4.Range("F12").Value = 4765
$ws4.Range("F13").Value = 580
$ws4.Range("F14").Value = 75
$ws4.Range("F15").Value = 26
$ws4.Range("F17").Value = 87
$ws4.Range("F18").Value = 934
$ws4.Range("F21").Value = 59
$ws4.Range("F22").Value = 5216
